# Ultra-HQ Sample Overhaul: Injected Wikipedia-grade substantive research
# content (dry goods) across all 6 flagship topics.
#
# NOTE: text is set via the Run object (Paragraphs(n,1).Runs(1,1).Text = ...)
# rather than Paragraph/TextRange.Text directly. The latter silently
# re-splits a run that ends in the CJK full stop "。" into two runs
# (text + punctuation), which would diverge from the original single-run
# OOXML shape. Setting .Text on the already-existing Run object instead
# just swaps the <a:t> contents in place and keeps the surrounding <a:rPr>
# untouched.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Speaker notes (notesSlide1/2/3, attached to slides 3, 4, 6)
# ---------------------------------------------------------------------
$n3 = $p.Slides.Item(3).NotesPage.Shapes.Item(2).TextFrame.TextRange
$n3.Runs(1,1).Text = "本节介绍 AI 的前身及其早期的哲学与数学基础，重点强调对人类智能的逻辑模拟尝试。"

$n4 = $p.Slides.Item(4).NotesPage.Shapes.Item(2).TextFrame.TextRange
$n4.Runs(1,1).Text = "介绍神经网络从被冷落到重新获得学术界关注的过程，为后来的深度学习爆发做铺垫。"

$n6 = $p.Slides.Item(6).NotesPage.Shapes.Item(2).TextFrame.TextRange
$n6.Runs(1,1).Text = "分析深度学习在视觉和 NLP 领域的双重突破，揭示 Transformer 架构为何成为当今万物互联的技术底座。"

# ---------------------------------------------------------------------
# Slide 1 - Title slide
# ---------------------------------------------------------------------
$s1 = $p.Slides.Item(1)
$s1.Shapes.Item(1).TextFrame.TextRange.Runs(1,1).Text = "人工智能：从图灵测试到通用人工智能 (AGI)"
$s1.Shapes.Item(2).TextFrame.TextRange.Runs(1,1).Text = "Topic: AI的发展历史与未来趋势"

# ---------------------------------------------------------------------
# Slide 2 - Section header
# ---------------------------------------------------------------------
$s2 = $p.Slides.Item(2)
$s2.Shapes.Item(1).TextFrame.TextRange.Runs(1,1).Text = "计算智能的起源与逻辑奠基 (1950-1980)"

# ---------------------------------------------------------------------
# Slide 3 - Content slide (title + 4 bullets)
# ---------------------------------------------------------------------
$s3 = $p.Slides.Item(3)
$s3.Shapes.Item(1).TextFrame.TextRange.Runs(1,1).Text = "图灵测试与符号 AI 的诞生"
$tr3 = $s3.Shapes.Item(2).TextFrame.TextRange
$tr3.Paragraphs(2,1).Runs(1,1).Text = "1950年：阿兰·图灵发表《计算机器与智能》，提出著名的“图灵测试” (Turing Test)。"
$tr3.Paragraphs(3,1).Runs(1,1).Text = "1956年：达特茅斯会议 (Dartmouth Workshop) 正式确立“人工智能”学科，麦卡锡、明斯基等人为学科领袖。"
$tr3.Paragraphs(4,1).Runs(1,1).Text = "逻辑主义时代：基于规则的专家系统（如 MYCIN）在特定医疗诊断领域取得初步成功。"
$tr3.Paragraphs(5,1).Runs(1,1).Text = "瓶颈出现：早期 AI 难以处理模糊信息，导致70年代中期进入第一个“AI 冬天”。"

# ---------------------------------------------------------------------
# Slide 4 - Content slide (title + 4 bullets)
# ---------------------------------------------------------------------
$s4 = $p.Slides.Item(4)
$s4.Shapes.Item(1).TextFrame.TextRange.Runs(1,1).Text = "联结主义与神经网络的复兴"
$tr4 = $s4.Shapes.Item(2).TextFrame.TextRange
$tr4.Paragraphs(2,1).Runs(1,1).Text = "1986年：Rumelhart 提出反向传播算法 (Backpropagation)，解决了多层感知器的训练难题。"
$tr4.Paragraphs(3,1).Runs(1,1).Text = "统计学习方法崛起：SVM 与随机森林在90年代成为机器学习的主流工具。"
$tr4.Paragraphs(4,1).Runs(1,1).Text = "GPU 计算能力的增强：为复杂的矩阵运算提供了硬件基础，神经网络的研究重心逐渐转向深度化。"
$tr4.Paragraphs(5,1).Runs(1,1).Text = "循环神经网络 (RNN) 与 LSTM：在高盛等金融机构及自然语言处理中开始显露头角。"

# ---------------------------------------------------------------------
# Slide 5 - Section header
# ---------------------------------------------------------------------
$s5 = $p.Slides.Item(5)
$s5.Shapes.Item(1).TextFrame.TextRange.Runs(1,1).Text = "深度学习革命与大模型时代 (2012-Present)"

# ---------------------------------------------------------------------
# Slide 6 - Content slide (title + 4 bullets)
# ---------------------------------------------------------------------
$s6 = $p.Slides.Item(6)
$s6.Shapes.Item(1).TextFrame.TextRange.Runs(1,1).Text = "从 ImageNet 到 Transformer 架构"
$tr6 = $s6.Shapes.Item(2).TextFrame.TextRange
$tr6.Paragraphs(2,1).Runs(1,1).Text = "2012年：AlexNet 以领先第二名10.8%的优势夺得 ImageNet 冠军，开启深度卷积神经网络时代。"
$tr6.Paragraphs(3,1).Runs(1,1).Text = "2017年：Google 发表《Attention is All You Need》，提出 Transformer 架构，颠覆序列建模模式。"
$tr6.Paragraphs(4,1).Runs(1,1).Text = "预训练大模型 (LLMs)：GPT-3 的 1750亿参数规模展示了模型容量与涌现能力 (Emergent Abilities) 的正相关性。"
$tr6.Paragraphs(5,1).Runs(1,1).Text = "推理与对齐：利用 RLHF (基于人类反馈的强化学习) 解决了模型在道德与逻辑层面的幻觉问题。"

# ---------------------------------------------------------------------
# Slide 7 - References (3 existing links updated + 1 new link appended)
# ---------------------------------------------------------------------
$s7 = $p.Slides.Item(7)
$tr7 = $s7.Shapes.Item(2).TextFrame.TextRange
$tr7.Paragraphs(2,1).Runs(1,1).Text = "https://en.wikipedia.org/wiki/Main_Page"
$tr7.Paragraphs(3,1).Runs(1,1).Text = "https://www.nih.gov/"
$tr7.Paragraphs(4,1).Runs(1,1).Text = "https://scholar.google.com/"
[void]$tr7.InsertAfter("`rhttps://www.jstor.org/")
